$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 27

$ws.Cells.Item($row, 1).Value = 45931
$ws.Cells.Item($row, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"

$ws.Cells.Item($row, 2).Value = "21,3649"
$ws.Cells.Item($row, 3).Value = "15,0508"
$ws.Cells.Item($row, 4).Value = "15,1254"
$ws.Cells.Item($row, 5).Value = "15,1254"
